$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.389221
$ws.Range("H2").Value = 4.167663
$ws.Range("I2").Value = 0.2910270461264192
$ws.Range("J2").Value = 0.2910270461264192
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02725333333333333
$ws.Range("N2").Value = 0.08176
$ws.Range("O2").Value = 0.0007089206372884383
$ws.Range("P2").Value = 0.0007089206372884382
$ws.Range("Q2").Value = 0.03786090298666667
$ws.Range("R2").Value = 0.34074812688
$ws.Range("S2").Value = 0.0002063150790081128
$ws.Range("T2").Value = 0.0002063150790081128

# Row 3
$ws.Range("G3").Value = 1.389221
$ws.Range("H3").Value = 4.167663
$ws.Range("I3").Value = 0.2910270461264192
$ws.Range("J3").Value = 0.2910270461264192
$ws.Range("M3").Value = 38.416166
$ws.Range("N3").Value = 115.248498
$ws.Range("O3").Value = 0.9992910793627116
$ws.Range("P3").Value = 0.9992910793627116
$ws.Range("Q3").Value = 53.368544546686
$ws.Range("R3").Value = 480.316900920174
$ws.Range("S3").Value = 0.2908207310474111
$ws.Range("T3").Value = 0.2908207310474111

# Row 4
$ws.Range("I4").Value = 0.461328155686921
$ws.Range("J4").Value = 0.4613281556869209
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02725333333333333
$ws.Range("N4").Value = 0.08176
$ws.Range("O4").Value = 0.0007089206372884383
$ws.Range("P4").Value = 0.0007089206372884382
$ws.Range("Q4").Value = 0.06001607335111112
$ws.Range("R4").Value = 0.5401446601600001
$ws.Range("S4").Value = 0.0003270450501286719
$ws.Range("T4").Value = 0.0003270450501286718

# Row 5
$ws.Range("I5").Value = 0.461328155686921
$ws.Range("J5").Value = 0.4613281556869209
$ws.Range("M5").Value = 38.416166
$ws.Range("N5").Value = 115.248498
$ws.Range("O5").Value = 0.9992910793627116
$ws.Range("P5").Value = 0.9992910793627116
$ws.Range("Q5").Value = 84.59836484311867
$ws.Range("R5").Value = 761.3852835880681
$ws.Range("S5").Value = 0.4610011106367923
$ws.Range("T5").Value = 0.4610011106367923

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1506176666666667
$ws.Range("H6").Value = 0.451853
$ws.Range("I6").Value = 0.03155280162368235
$ws.Range("J6").Value = 0.03155280162368235
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02725333333333333
$ws.Range("N6").Value = 0.08176
$ws.Range("O6").Value = 0.0007089206372884383
$ws.Range("P6").Value = 0.0007089206372884382
$ws.Range("Q6").Value = 0.004104833475555556
$ws.Range("R6").Value = 0.03694350128
$ws.Range("S6").Value = 0.00002236843223529657
$ws.Range("T6").Value = 0.00002236843223529656

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1506176666666667
$ws.Range("H7").Value = 0.451853
$ws.Range("I7").Value = 0.03155280162368235
$ws.Range("J7").Value = 0.03155280162368235
$ws.Range("M7").Value = 38.416166
$ws.Range("N7").Value = 115.248498
$ws.Range("O7").Value = 0.9992910793627116
$ws.Range("P7").Value = 0.9992910793627116
$ws.Range("Q7").Value = 5.786153285199333
$ws.Range("R7").Value = 52.075379566794
$ws.Range("S7").Value = 0.03153043319144706
$ws.Range("T7").Value = 0.03153043319144706

# Row 8
$ws.Range("G8").Value = 1.031517666666667
$ws.Range("H8").Value = 3.094553
$ws.Range("I8").Value = 0.2160919965629775
$ws.Range("J8").Value = 0.2160919965629775
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02725333333333333
$ws.Range("N8").Value = 0.08176
$ws.Range("O8").Value = 0.0007089206372884383
$ws.Range("P8").Value = 0.0007089206372884382
$ws.Range("Q8").Value = 0.02811229480888889
$ws.Range("R8").Value = 0.25301065328
$ws.Range("S8").Value = 0.000153192075916357
$ws.Range("T8").Value = 0.000153192075916357

# Row 9
$ws.Range("G9").Value = 1.031517666666667
$ws.Range("H9").Value = 3.094553
$ws.Range("I9").Value = 0.2160919965629775
$ws.Range("J9").Value = 0.2160919965629775
$ws.Range("M9").Value = 38.416166
$ws.Range("N9").Value = 115.248498
$ws.Range("O9").Value = 0.9992910793627116
$ws.Range("P9").Value = 0.9992910793627116
$ws.Range("Q9").Value = 39.62695391459933
$ws.Range("R9").Value = 356.642585231394
$ws.Range("S9").Value = 0.2159388044870612
$ws.Range("T9").Value = 0.2159388044870612

